# NIT-9012107125.xlsx update:
#  - Update "VALOR MORA" total (E11) and "Cant. Periodos" count (F13)
#  - Add a 4th "Periodo Mora" row (2508) for the existing worker, reusing the
#    same layout/borders as the existing 3 rows (2505/2506/2507), which pushes
#    the bottom-border styling from the old last data row down to the new one
#  - Renumber the existing period rows (2507/2506/2505 -> 2505/2506/2507/2508)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -------------------------------------------------
$ws.Range("E11").Value = 227760   # VALOR MORA (was 170820)
$ws.Range("F13").Value = 4        # Cant. Periodos (was 3)

# --- Insert a new "Periodo Mora" row below the current last row (18) -------
$ws.Rows("19:19").Insert()

# The new row should look like the (still unmodified) old row 18 - i.e. the
# "last row" border treatment - so copy its formatting down first.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats

# Row 18 stops being the last row, so it now takes the "middle row" border
# treatment that rows 16/17 already use.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Fill in the new row's data (same worker, new period) -------------------
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1050964669"
$ws.Range("D19").Value = "ALDO ENRIQUE GONZALEZ DIAZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- Renumber the periods shown in rows 16-18 -------------------------------
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
